$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.836.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.129.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.16%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.16%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'600.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.08%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'140.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.16%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.130.92"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.99%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.36%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.91%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.33%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.21%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'34.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.53%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.646.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.13%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +2.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'63.872.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.33%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.139.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.86%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -1.89%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'482.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.80%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.706"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.15%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.05%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'87.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.14%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'13.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.44%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.02%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -2.66%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'8.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -7.03%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'6.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.47%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.11%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'27.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.52%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -7.05%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.23%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.65%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.84%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.08%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'52.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.24%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0₃0729"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -8.00%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0396"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.83%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -9.84%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'428.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -7.18%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.75%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.56%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.898.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.28%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -3.64%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -7.35%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E48").Value = "'  -3.22%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.26%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'25.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.92%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'120.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.60%  "
$ws.Range("E51").Style = "Normal"
